$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.716.33'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '2.372.39'
$ws.Range("E3").Value = '  -3.24%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.84%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -10.48%  '

$ws.Range("D9").Value = '2.369.86'
$ws.Range("E9").Value = '  -3.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.105'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.89%  '

$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.343'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.33%  '

$ws.Range("D15").Value = '2.794.61'
$ws.Range("E15").Value = '  -3.55%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000162'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '60.488.63'
$ws.Range("E17").Value = '  -0.37%  '

$ws.Range("D18").Value = '2.364.05'
$ws.Range("E18").Value = '  -3.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '314.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").Value = '2.476.36'
$ws.Range("E27").Value = '  -3.98%  '

$ws.Range("D28").Value = '0.0₃0923'
$ws.Range("E28").Value = '  -5.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '519.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.35%  '

$ws.Range("E33").Value = '  -3.54%  '

$ws.Range("E34").Value = '  -3.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.07%  '

$ws.Range("E36").Value = '  -0.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.44'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.374'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.80%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("E45").Value = '  -2.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0514'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.574'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0906'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.23%  '

